$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.238.24"
Set-TextValue $ws.Range("E2") "  +1.79%  "
Set-TextValue $ws.Range("D3") "1.890.42"
Set-TextValue $ws.Range("E3") "  -1.30%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "323.44"
Set-TextValue $ws.Range("E5") "  +1.75%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.01%  "
Set-TextValue $ws.Range("D7") "0.5176"
Set-TextValue $ws.Range("E7") "  -0.16%  "
Set-TextValue $ws.Range("D8") "0.4011"
Set-TextValue $ws.Range("E8") "  +1.01%  "
Set-TextValue $ws.Range("D9") "0.08391"
Set-TextValue $ws.Range("E9") "  -1.48%  "
Set-TextValue $ws.Range("D10") "42.68"
Set-TextValue $ws.Range("E10") "  -0.09%  "
Set-TextValue $ws.Range("D11") "1.112"
Set-TextValue $ws.Range("E11") "  -0.92%  "
Set-TextValue $ws.Range("D12") "23.05"
Set-TextValue $ws.Range("E12") "  +10.16%  "
Set-TextValue $ws.Range("D13") "6.421"
Set-TextValue $ws.Range("E13") "  +1.96%  "
Set-TextValue $ws.Range("D14") "1.894.07"
Set-TextValue $ws.Range("E14") "  -0.76%  "
Set-TextValue $ws.Range("D15") "7.305"
Set-TextValue $ws.Range("E15") "  -0.68%  "
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  +0.01%  "
Set-TextValue $ws.Range("D17") "94.10"
Set-TextValue $ws.Range("E17") "  +0.03%  "
Set-TextValue $ws.Range("D18") "0.00001107"
Set-TextValue $ws.Range("E18") "  -0.71%  "
Set-TextValue $ws.Range("D19") "0.06638"
Set-TextValue $ws.Range("E19") "  -1.63%  "
Set-TextValue $ws.Range("D20") "18.19"
Set-TextValue $ws.Range("E20") "  +1.27%  "
Set-TextValue $ws.Range("E21") "  -0.01%  "
Set-TextValue $ws.Range("D22") "5.937"
Set-TextValue $ws.Range("E22") "  -1.67%  "
Set-TextValue $ws.Range("D23") "30.233.84"
Set-TextValue $ws.Range("E23") "  +1.76%  "
Set-TextValue $ws.Range("D24") "11.26"
Set-TextValue $ws.Range("E24") "  +0.51%  "
Set-TextValue $ws.Range("D25") "2.229"
Set-TextValue $ws.Range("E25") "  +0.80%  "
Set-TextValue $ws.Range("D26") "2.114.88"
Set-TextValue $ws.Range("E26") "  -0.58%  "
Set-TextValue $ws.Range("D27") "21.55"
Set-TextValue $ws.Range("E27") "  +2.17%  "
Set-TextValue $ws.Range("D28") "161.78"
Set-TextValue $ws.Range("E28") "  +1.51%  "
Set-TextValue $ws.Range("D29") "2.324"
Set-TextValue $ws.Range("E29") "  -5.79%  "
Set-TextValue $ws.Range("D30") "128.91"
Set-TextValue $ws.Range("E30") "  +0.09%  "
Set-TextValue $ws.Range("D31") "1.085"
Set-TextValue $ws.Range("E31") "  -0.08%  "
Set-TextValue $ws.Range("D32") "0.1052"
Set-TextValue $ws.Range("E32") "  -0.68%  "
Set-TextValue $ws.Range("D33") "6.081"
Set-TextValue $ws.Range("E33") "  -2.01%  "
Set-TextValue $ws.Range("D34") "3.747"
Set-TextValue $ws.Range("E34") "  +1.77%  "
Set-TextValue $ws.Range("D35") "0.02486"
Set-TextValue $ws.Range("E35") "  -0.50%  "
Set-TextValue $ws.Range("D36") "0.06526"
Set-TextValue $ws.Range("E36") "  -1.69%  "
Set-TextValue $ws.Range("D37") "5.329"
Set-TextValue $ws.Range("E37") "  +2.20%  "
Set-TextValue $ws.Range("D38") "0.2191"
Set-TextValue $ws.Range("E38") "  -0.71%  "
Set-TextValue $ws.Range("D39") "1.217"
Set-TextValue $ws.Range("E39") "  -2.70%  "
Set-TextValue $ws.Range("D40") "8.803"
Set-TextValue $ws.Range("E40") "  -3.43%  "
Set-TextValue $ws.Range("D41") "11.76"
Set-TextValue $ws.Range("E41") "  +3.39%  "
Set-TextValue $ws.Range("D42") "0.6479"
Set-TextValue $ws.Range("E42") "  -1.02%  "
Set-TextValue $ws.Range("D43") "1.227"
Set-TextValue $ws.Range("E43") "  -1.18%  "
Set-TextValue $ws.Range("D44") "0.6071"
Set-TextValue $ws.Range("E44") "  -0.95%  "
Set-TextValue $ws.Range("E45") "  -0.36%  "
Set-TextValue $ws.Range("D46") "3.689"
Set-TextValue $ws.Range("E46") "  -0.04%  "
Set-TextValue $ws.Range("D47") "2.049"
Set-TextValue $ws.Range("E47") "  -0.82%  "
Set-TextValue $ws.Range("D48") "1.235"
Set-TextValue $ws.Range("E48") "  -0.38%  "
Set-TextValue $ws.Range("D49") "124.45"
Set-TextValue $ws.Range("E49") "  -0.30%  "
Set-TextValue $ws.Range("D50") "1.159"
Set-TextValue $ws.Range("E50") "  -2.34%  "
Set-TextValue $ws.Range("D51") "78.89"
Set-TextValue $ws.Range("E51") "  +0.51%  "
